$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '50.887.02'
$ws.Range('E2').Value = '  -0.77%  '

$ws.Range('D3').Value = '2.926.09'
$ws.Range('E3').Value = '  -1.27%  '

$ws.Range('E4').Value = '  +0.17%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '372.90'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -1.90%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '100.84'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -4.09%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.534'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -1.32%  '

$ws.Range('E8').Value = '  -0.03%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.581'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -1.97%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '35.90'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -3.66%  '

$ws.Range('E11').Value = '  -0.60%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.0840'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -0.74%  '

$ws.Range('D13').Value = '3.392.78'
$ws.Range('E13').Value = '  -1.09%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '17.84'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -3.03%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.37'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -2.18%  '

$ws.Range('B16').Value = 'Uniswap'
$ws.Range('C16').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '11.34'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +53.92%  '

$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '2.955.06'
$ws.Range('E17').Value = '  -0.09%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.968'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +0.49%  '

$ws.Range('D19').Value = '50.817.71'
$ws.Range('E19').Value = '  -0.85%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '3.12'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -6.13%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '12.42'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -3.59%  '

$ws.Range('D22').Value = '0.0₃0951'
$ws.Range('E22').Value = '  -1.14%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '263.34'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +0.94%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '68.32'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -1.59%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '3.09'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +9.64%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '7.97'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +4.09%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '7.54'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +1.94%  '

$ws.Range('B28').Value = 'Dai'
$ws.Range('C28').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.999'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -0.11%  '

$ws.Range('B29').Value = 'Kaspa'
$ws.Range('C29').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.166'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -2.95%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '25.47'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -1.36%  '

$ws.Range('E31').Value = '  -2.33%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '9.92'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +0.28%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '50.57'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -1.29%  '

$ws.Range('E34').Value = '  -2.94%  '

$ws.Range('B35').Value = 'VeChain'
$ws.Range('C35').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.0438'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -1.98%  '

$ws.Range('B36').Value = 'InjectiveProtocol'
$ws.Range('C36').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '32.67'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -7.00%  '

$ws.Range('E37').Value = '  -0.01%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '3.10'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +1.61%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.115'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -0.62%  '

$ws.Range('B40').Value = 'Celestia'
$ws.Range('C40').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '16.12'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -6.47%  '

$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.47'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -4.27%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.77'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -3.96%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '120.24'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -3.33%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '20.90'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -4.61%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.276'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -4.45%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.03'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -1.94%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '3.27'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +0.93%  '

$ws.Range('E48').Value = '  -3.28%  '

$ws.Range('D49').Value = '1.987.45'
$ws.Range('E49').Value = '  -2.60%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0331'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -3.62%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.28'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -0.22%  '
